$wb = $excel.ActiveWorkbook

# --- Sheet "Student Points": clear the "No exercise screen-shots" note
#     and record full credit for the exercise ---
$ws2 = $wb.Worksheets.Item("Student Points")
$ws2.Range("D7").ClearContents()
$ws2.Range("C7").Value = 10
[void]$ws2.Range("A4:C23").Select()

# --- Sheet "Lab3Rubric_CS295N": rubric point redistribution in the
#     "Unit tests" section (rows 12-13). This sheet stays active/selected,
#     so touch it last. ---
$ws1 = $wb.Worksheets.Item("Lab3Rubric_CS295N")
$ws1.Range("B12").Value = 5
$ws1.Range("C12").Value = 5
$ws1.Range("B13").Value = 5
$ws1.Range("C13").Value = 5
[void]$ws1.Range("E6").Select()
